$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.020.29"
$ws.Range("E2").Value = '  -1.66%  '

$ws.Range("D3").Value = "'2.943.95"
$ws.Range("E3").Value = '  -2.12%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = "'586.79"
$ws.Range("E5").Value = '  -1.73%  '

$ws.Range("D6").Value = "'146.86"
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").Value = "'2.926.27"
$ws.Range("E8").Value = '  -2.72%  '

$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = '  -2.63%  '

$ws.Range("D10").Value = "'6.77"
$ws.Range("E10").Value = '  +7.76%  '

$ws.Range("D11").Value = "'0.145"
$ws.Range("E11").Value = '  -2.82%  '

$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = '  -1.33%  '

$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = '  -1.98%  '

$ws.Range("D14").Value = "'34.54"
$ws.Range("E14").Value = '  +0.44%  '

$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").Value = "'3.432.74"
$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").Value = "'6.87"
$ws.Range("E17").Value = '  -1.21%  '

$ws.Range("D18").Value = "'61.148.77"
$ws.Range("E18").Value = '  -1.32%  '

$ws.Range("D19").Value = "'2.944.87"
$ws.Range("E19").Value = '  -2.23%  '

$ws.Range("D20").Value = "'430.49"
$ws.Range("E20").Value = '  -4.23%  '

$ws.Range("D21").Value = "'13.84"
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").Value = "'0.674"
$ws.Range("E22").Value = '  -1.84%  '

$ws.Range("D23").Value = "'7.27"
$ws.Range("E23").Value = '  -1.54%  '

$ws.Range("D24").Value = "'80.56"
$ws.Range("E24").Value = '  -1.52%  '

$ws.Range("D25").Value = "'10.92"
$ws.Range("E25").Value = '  -1.55%  '

$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").Value = "'11.92"
$ws.Range("E27").Value = '  -2.18%  '

$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = '  +1.90%  '

$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("E31").Value = '  +5.31%  '

$ws.Range("E32").Value = '  -2.68%  '

$ws.Range("D33").Value = "'26.93"
$ws.Range("E33").Value = '  -1.41%  '

$ws.Range("E34").Value = '  -2.87%  '

$ws.Range("D35").Value = "'0.0₃0831"
$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("E36").Value = '  -1.38%  '

$ws.Range("D37").Value = "'5.71"
$ws.Range("E37").Value = '  -1.63%  '

$ws.Range("D38").Value = "'49.96"
$ws.Range("E38").Value = '  -1.11%  '

$ws.Range("D39").Value = "'2.99"
$ws.Range("E39").Value = '  +2.74%  '

$ws.Range("D40").Value = "'0.126"
$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("D41").Value = "'2.03"
$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("D42").Value = "'8.77"
$ws.Range("E42").Value = '  -2.58%  '

$ws.Range("D43").Value = "'0.291"
$ws.Range("E43").Value = '  +6.28%  '

$ws.Range("D44").Value = "'41.58"
$ws.Range("E44").Value = '  +2.03%  '

$ws.Range("D45").Value = "'0.0348"
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("D46").Value = "'372.76"
$ws.Range("E46").Value = '  -6.54%  '

$ws.Range("D47").Value = "'2.660.13"
$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("D48").Value = "'133.29"
$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("D49").Value = "'25.42"
$ws.Range("E49").Value = '  +7.44%  '

$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = '  -0.87%  '
